# "All \items with a capital, make sure issues tbl is consistnt with source"
#
# Refreshes the window view state and corrects six source numbers on the
# "Data for all bm (2)" sheet's issues table so it matches upstream data.
# The sheet's dependent formulas (the "Casts per 100 LOC" ratio row, and
# the mirrored rows in the "TABLE FOR THESIS" block further down) are
# driven off these cells, so they recalculate automatically.

$wb = $excel.ActiveWorkbook

# Best-effort refresh of the saved window geometry (xWindow/yWindow/
# windowWidth/windowHeight in the workbook's bookViews).
$win = $excel.ActiveWindow
$win.Top    = -20
$win.Left   = 0
$win.Width  = 33600
$win.Height = 20540

# The edited sheet.
$ws = $wb.Worksheets.Item("Data for all bm (2)")
$ws.Activate()

# Row 9 ("Better lang. support for shorts and bytes" / "Casts"): counts.
$ws.Range("K9").Value = 10
$ws.Range("O9").Value = 64

# Row 10 (the corresponding "LOC" row): counts.
$ws.Range("J10").Value = 44
$ws.Range("K10").Value = 77
$ws.Range("N10").Value = 51
$ws.Range("O10").Value = 266

# Move the sheet's active selection to match the latest edited cell.
$ws.Range("O11").Select()

$excel.Calculate()
